# Update simulation results: new "min load val" (column B) and
# "min load time" (column C) outputs for rows 2-11.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New minimum load values (column B)
$bValues = @(
    0.9559800000000001,
    0.95174,
    0.95814,
    0.96065,
    0.9545400000000001,
    0.96583,
    0.96099,
    0.96236,
    0.95382,
    0.96184
)

# New minimum load times (column C), as text
$cValues = @(
    "18:40:00",
    "18:27:00",
    "18:33:00",
    "18:45:00",
    "18:49:00",
    "18:52:00",
    "18:37:00",
    "18:48:00",
    "18:33:00",
    "18:28:00"
)

for ($i = 0; $i -lt $bValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 2).Value = $bValues[$i]
    $ws.Cells.Item($row, 3).Value = $cValues[$i]
}
